# ERSP essay - second draft rewrite
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the six body paragraphs (index 3..8) with the revised text, and
#    add a first-line indent to each of them.
# ---------------------------------------------------------------------------

$newTexts = @(
  "As an international student from a modest socioeconomic background, I have faced several challenges – notably, that of growing up in a small hometown with limited resources and opportunities. However, it has taught me to make the most out of every opportunity and approach challenges with resilience, and I believe that is an invaluable perspective I can bring to the Early Research Scholars Program. I relate to others coming from similarly modest or low socioeconomic backgrounds, and by voicing out my ideas and contributing to the group, I would promote and inspire their contribution. This reminds me of a quote from the film `Coach Carter', which goes `As we let our own light shine, we unconsciously give other people permission to do the same.'",
  "I understand that diversity is not just about representation but about fostering an inclusive space where everyone and their contributions are acknowledged and valued. I would actively strive to be a good team player, respecting people's identities and fostering an atmosphere where ideas can flow freely. For instance, I would make it a priority to ask for and use people's desired pronouns, ensuring that non-binary members feel validated. Such simple acts of respect can go a long way in creating a sense of belonging. ",
  "Additionally, I would work to amplify the voices of any students from historically disadvantaged communities, acknowledging their challenges and celebrating the unique perspectives they bring to the table. By encouraging open dialogue and collaboration, I hope to build an environment where everyone feels comfortable sharing their ideas. Even if it is just a germ of an idea that isn't yet well-developed, I would encourage them to speak out their mind. There's a good chance that there is utility in that idea, and, if there is, we could always work together to improvise on that idea. Encouraging open dialogue and making everyone feel comfortable sharing their ideas would be my priority.",
  "For those who may have experienced trauma, from, say, being raised in violent nations or other difficult life circumstances, I would be just as considerate of their life experiences. Plus, I believe in encouragement and positive reinforcement. So, good ideas would be appreciated, and not-so-helpful contributions would still be appreciated! The effort is what matters. The goal would be to foster a space free of judgement, where everyone feels safe and included – whether it is in group interactions or research activities.",
  "I also understand how important it can be to address accessibility concerns for people with disabilities. I would attempt to ensure materials are usable by everyone, and adjust my communication and other aspects of group interactions to accommodate their needs. Creating a collaborative environment that highlights everyone's strengths, rather than their limitations, is an important step in building a truly inclusive community. ",
  "Ultimately, I believe that the strength of ERSP lies in the diversity of its participants. By listening and learning form their inputs, and supporting their identities, I hope to participate in a community where everyone feels empowered to achieve and contribute. I'm committed to helping ERSP remain a space where students from all walks of life feel inspired and valued."
)

for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $idx = 3 + $i
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $body = $d.Range($r.Start, $r.End - 1)
    $body.Text = $newTexts[$i]
    $p.Range.ParagraphFormat.FirstLineIndent = 36
}

# ---------------------------------------------------------------------------
# 2. Apply the strikethrough formatting inside the new "Additionally" paragraph
#    (index 5) on the sentence that got struck through in the revision.
# ---------------------------------------------------------------------------

$strikeRng = $d.Content
$strikeRng.Find.ClearFormatting()
$found = $strikeRng.Find.Execute("By encouraging open dialogue and collaboration, I hope to build an environment where everyone feels comfortable sharing their ideas", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $strikeRng.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# 3. Rebuild the comments: drop all five legacy comments and re-add the three
#    that survive in the revision (in left-to-right order so the new `w:id`
#    numbering comes out as 0, 1, 2).
# ---------------------------------------------------------------------------

for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# Comment 0: attached to "communities" in the "Additionally" paragraph.
$rng0 = $d.Content
$rng0.Find.ClearFormatting()
$rng0.Find.Execute("communities", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c0 = $d.Comments.Add($rng0, "Covers `"races and ethnicities.'")
$c0.Author = "Prateek Basavaraj"
$c0.Initial = "PB"

# Comment 1: attached to "encouragement" in the "For those who may have
# experienced trauma" paragraph.
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("encouragement", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c1 = $d.Comments.Add($rng1, "`'encouragement' - Superfluous word? (positive reinforcement already describes the same thing)")
$c1.Author = "Prateek Basavaraj"
$c1.Initial = "PB"

# Comment 2: attached to "achieve and contribute" in the final "Ultimately" paragraph.
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("achieve and contribute", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c2 = $d.Comments.Add($rng2, "superfluous")
$c2.Author = "Prateek Basavaraj"
$c2.Initial = "PB"

# ---------------------------------------------------------------------------
# 4. Tidy up the bullet list item about presentation/communication skills:
#    merge the three runs that used to be split by proofing marks into one.
# ---------------------------------------------------------------------------

$skillsRng = $d.Content
$skillsRng.Find.ClearFormatting()
$skillsRng.Find.Execute("skills(", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullSkillsRng = $d.Content
$fullSkillsRng.Find.ClearFormatting()
$fullSkillsRng.Find.Execute("can boost presentation skills and communication skills(team lingo, etc.) – invaluable soft skills.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $fullSkillsRng.Find.Found) {
    $anchor = $d.Content
    $anchor.Find.ClearFormatting()
    $anchor.Find.Execute("can boost presentation skills and communication ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
Write-Host "done"
